# Insert a new data row at row 72 (pushing existing rows 72..187 down to 73..188)
# and populate it with the new "Perejil" price record described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

$ws.Cells.Item(72, 1).Value  = 4
$ws.Cells.Item(72, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value  = "Los Lagos"
$ws.Cells.Item(72, 4).Value  = 44540
$ws.Cells.Item(72, 5).Value  = 10
$ws.Cells.Item(72, 6).Value  = 100112044
$ws.Cells.Item(72, 7).Value  = "Perejil"
$ws.Cells.Item(72, 8).Value  = "Sin especificar"
$ws.Cells.Item(72, 9).Value  = "Primera"
$ws.Cells.Item(72, 10).Value = 180
$ws.Cells.Item(72, 11).Value = 5000
$ws.Cells.Item(72, 12).Value = 5000
$ws.Cells.Item(72, 13).Value = 5000
$ws.Cells.Item(72, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 1667
$ws.Cells.Item(72, 17).Value = 3
$ws.Cells.Item(72, 18).Value = "Hortaliza"
